$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "3.575558"
$ws.Range("H2").Value = "10.726674"
$ws.Range("I2").Value = "0.025194653521236"
$ws.Range("J2").Value = "0.02519465352123599"
$ws.Range("M2").Value = "8.820647333333334"
$ws.Range("N2").Value = "26.461942"
$ws.Range("O2").Value = "0.06415146660411865"
$ws.Range("P2").Value = "0.06415146660411865"
$ws.Range("Q2").Value = "31.53873613787867"
$ws.Range("R2").Value = "283.848625240908"
$ws.Range("S2").Value = "0.001616273973969911"
$ws.Range("T2").Value = "0.001616273973969911"
$ws.Range("G3").Value = "3.575558"
$ws.Range("H3").Value = "10.726674"
$ws.Range("I3").Value = "0.025194653521236"
$ws.Range("J3").Value = "0.02519465352123599"
$ws.Range("O3").Value = "0.3979101621202897"
$ws.Range("P3").Value = "0.3979101621202898"
$ws.Range("Q3").Value = "195.62426666153"
$ws.Range("R3").Value = "1760.61839995377"
$ws.Range("S3").Value = "0.01002520866719954"
$ws.Range("T3").Value = "0.01002520866719954"
$ws.Range("G4").Value = "3.575558"
$ws.Range("H4").Value = "10.726674"
$ws.Range("I4").Value = "0.025194653521236"
$ws.Range("J4").Value = "0.02519465352123599"
$ws.Range("M4").Value = "21.90816333333333"
$ws.Range("N4").Value = "65.72449"
$ws.Range("O4").Value = "0.1593353362087987"
$ws.Range("P4").Value = "0.1593353362087987"
$ws.Range("Q4").Value = "78.33390867180665"
$ws.Range("R4").Value = "705.00517804626"
$ws.Range("S4").Value = "0.004014398589470331"
$ws.Range("T4").Value = "0.00401439858947033"
$ws.Range("G5").Value = "3.575558"
$ws.Range("H5").Value = "10.726674"
$ws.Range("I5").Value = "0.025194653521236"
$ws.Range("J5").Value = "0.02519465352123599"
$ws.Range("M5").Value = "52.056859"
$ws.Range("N5").Value = "156.170577"
$ws.Range("O5").Value = "0.3786030350667928"
$ws.Range("P5").Value = "0.3786030350667929"
$ws.Range("Q5").Value = "186.132318652322"
$ws.Range("R5").Value = "1675.190867870898"
$ws.Range("S5").Value = "0.009538772290596208"
$ws.Range("T5").Value = "0.009538772290596208"
$ws.Range("I6").Value = "0.7460690747908298"
$ws.Range("J6").Value = "0.7460690747908298"
$ws.Range("M6").Value = "8.820647333333334"
$ws.Range("N6").Value = "26.461942"
$ws.Range("O6").Value = "0.06415146660411865"
$ws.Range("P6").Value = "0.06415146660411865"
$ws.Range("Q6").Value = "933.9313069189972"
$ws.Range("R6").Value = "8405.381762270974"
$ws.Range("S6").Value = "0.04786142533580962"
$ws.Range("T6").Value = "0.04786142533580962"
$ws.Range("I7").Value = "0.7460690747908298"
$ws.Range("J7").Value = "0.7460690747908298"
$ws.Range("O7").Value = "0.3979101621202897"
$ws.Range("P7").Value = "0.3979101621202898"
$ws.Range("S7").Value = "0.2968684665029537"
$ws.Range("T7").Value = "0.2968684665029537"
$ws.Range("I8").Value = "0.7460690747908298"
$ws.Range("J8").Value = "0.7460690747908298"
$ws.Range("M8").Value = "21.90816333333333"
$ws.Range("N8").Value = "65.72449"
$ws.Range("O8").Value = "0.1593353362087987"
$ws.Range("P8").Value = "0.1593353362087987"
$ws.Range("Q8").Value = "2319.639232913615"
$ws.Range("R8").Value = "20876.75309622253"
$ws.Range("S8").Value = "0.1188751668667842"
$ws.Range("T8").Value = "0.1188751668667842"
$ws.Range("I9").Value = "0.7460690747908298"
$ws.Range("J9").Value = "0.7460690747908298"
$ws.Range("M9").Value = "52.056859"
$ws.Range("N9").Value = "156.170577"
$ws.Range("O9").Value = "0.3786030350667928"
$ws.Range("P9").Value = "0.3786030350667929"
$ws.Range("Q9").Value = "5511.787119777674"
$ws.Range("R9").Value = "49606.08407799906"
$ws.Range("S9").Value = "0.2824640160852822"
$ws.Range("T9").Value = "0.2824640160852823"
$ws.Range("G10").Value = "32.36130266666667"
$ws.Range("H10").Value = "97.08390800000001"
$ws.Range("I10").Value = "0.2280292497513723"
$ws.Range("J10").Value = "0.2280292497513723"
$ws.Range("M10").Value = "8.820647333333334"
$ws.Range("N10").Value = "26.461942"
$ws.Range("O10").Value = "0.06415146660411865"
$ws.Range("P10").Value = "0.06415146660411865"
$ws.Range("Q10").Value = "285.4476380699263"
$ws.Range("R10").Value = "2569.028742629336"
$ws.Range("S10").Value = "0.01462841080018739"
$ws.Range("T10").Value = "0.01462841080018739"
$ws.Range("G11").Value = "32.36130266666667"
$ws.Range("H11").Value = "97.08390800000001"
$ws.Range("I11").Value = "0.2280292497513723"
$ws.Range("J11").Value = "0.2280292497513723"
$ws.Range("O11").Value = "0.3979101621202897"
$ws.Range("P11").Value = "0.3979101621202898"
$ws.Range("Q11").Value = "1770.536543492927"
$ws.Range("R11").Value = "15934.82889143634"
$ws.Range("S11").Value = "0.09073515573673659"
$ws.Range("T11").Value = "0.0907351557367366"
$ws.Range("G12").Value = "32.36130266666667"
$ws.Range("H12").Value = "97.08390800000001"
$ws.Range("I12").Value = "0.2280292497513723"
$ws.Range("J12").Value = "0.2280292497513723"
$ws.Range("M12").Value = "21.90816333333333"
$ws.Range("N12").Value = "65.72449"
$ws.Range("O12").Value = "0.1593353362087987"
$ws.Range("P12").Value = "0.1593353362087987"
$ws.Range("Q12").Value = "708.976704500769"
$ws.Range("R12").Value = "6380.790340506921"
$ws.Range("S12").Value = "0.03633311717457503"
$ws.Range("T12").Value = "0.03633311717457503"
$ws.Range("G13").Value = "32.36130266666667"
$ws.Range("H13").Value = "97.08390800000001"
$ws.Range("I13").Value = "0.2280292497513723"
$ws.Range("J13").Value = "0.2280292497513723"
$ws.Range("M13").Value = "52.056859"
$ws.Range("N13").Value = "156.170577"
$ws.Range("O13").Value = "0.3786030350667928"
$ws.Range("P13").Value = "0.3786030350667929"
$ws.Range("Q13").Value = "1684.627769974991"
$ws.Range("R13").Value = "15161.64992977492"
$ws.Range("S13").Value = "0.08633256603987326"
$ws.Range("T13").Value = "0.08633256603987327"
$ws.Range("G14").Value = "0.1003386666666667"
$ws.Range("H14").Value = "0.301016"
$ws.Range("I14").Value = "0.000707021936561918"
$ws.Range("J14").Value = "0.0007070219365619179"
$ws.Range("M14").Value = "8.820647333333334"
$ws.Range("N14").Value = "26.461942"
$ws.Range("O14").Value = "0.06415146660411865"
$ws.Range("P14").Value = "0.06415146660411865"
$ws.Range("Q14").Value = "0.8850519925635556"
$ws.Range("R14").Value = "7.965467933072"
$ws.Range("S14").Value = "0.00004535649415173118"
$ws.Range("T14").Value = "0.00004535649415173117"
$ws.Range("G15").Value = "0.1003386666666667"
$ws.Range("H15").Value = "0.301016"
$ws.Range("I15").Value = "0.000707021936561918"
$ws.Range("J15").Value = "0.0007070219365619179"
$ws.Range("O15").Value = "0.3979101621202897"
$ws.Range("P15").Value = "0.3979101621202898"
$ws.Range("Q15").Value = "5.489682473186667"
$ws.Range("R15").Value = "49.40714225868"
$ws.Range("S15").Value = "0.000281331213399954"
$ws.Range("T15").Value = "0.000281331213399954"
$ws.Range("G16").Value = "0.1003386666666667"
$ws.Range("H16").Value = "0.301016"
$ws.Range("I16").Value = "0.000707021936561918"
$ws.Range("J16").Value = "0.0007070219365619179"
$ws.Range("M16").Value = "21.90816333333333"
$ws.Range("N16").Value = "65.72449"
$ws.Range("O16").Value = "0.1593353362087987"
$ws.Range("P16").Value = "0.1593353362087987"
$ws.Range("Q16").Value = "2.198235897982222"
$ws.Range("R16").Value = "19.78412308184"
$ws.Range("S16").Value = "0.0001126535779690891"
$ws.Range("T16").Value = "0.0001126535779690891"
$ws.Range("G17").Value = "0.1003386666666667"
$ws.Range("H17").Value = "0.301016"
$ws.Range("I17").Value = "0.000707021936561918"
$ws.Range("J17").Value = "0.0007070219365619179"
$ws.Range("M17").Value = "52.056859"
$ws.Range("N17").Value = "156.170577"
$ws.Range("O17").Value = "0.3786030350667928"
$ws.Range("P17").Value = "0.3786030350667929"
$ws.Range("Q17").Value = "5.223315822914667"
$ws.Range("R17").Value = "47.009842406232"
$ws.Range("S17").Value = "0.0002676806510411436"
$ws.Range("T17").Value = "0.0002676806510411436"
